$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update F2 value
$ws.Range("F2").Value = 680380

# Clear F3 (value removed entirely)
$ws.Range("F3").ClearContents()

# F4 and F5 are formulas (F2-F3 and F4/F2), they recalc automatically.
# Update the selected cell to H12 (matches the <selection> change in the diff)
$ws.Range("H12").Select()
